# Added Short Screenshot to file.
# Slide 2 has a textbox ("TextBox 9", id=10) with two stacked labels "Sc" / "Sh"
# describing the "Screenshot" / "Shift" combo key. This edit:
#   1. Turns it into "ScSh" / "Cp" (Screenshot+Shift -> Copy), shrinking it and
#      moving it slightly.
#   2. Duplicates it into a new textbox ("TextBox 6", id=7) labelled
#      "ScSh" / "File" (Screenshot+Shift -> File) placed to its left.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# Locate the existing "Sc"/"Sh" textbox by name.
$orig = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.Name -eq "TextBox 9") {
        $orig = $candidate
        break
    }
}

# Duplicate it first so the clone inherits all of the original's formatting
# (body/run properties, language, color, autofit, etc.) before we touch it.
$dup = $orig.Duplicate()
$dup.Name = "TextBox 6"

# --- Update the original shape: "Sc"/"Sh" -> "ScSh"/"Cp" ---
# (Text/font must be set before size, since the textbox has spAutoFit and
# PowerPoint recalculates Height from the text right after it changes.)
$orig.TextFrame.TextRange.Text = "ZZZZ"
$orig.TextFrame.TextRange.Text = "ScSh`rCp"
$orig.TextFrame.TextRange.Font.Size = 15

$orig.Left = 256.27906801811025
$orig.Top = 244.4292984385827
$orig.Width = 48.52291298582677
$orig.Height = 43.62188916377953

# --- Configure the duplicated shape: "ScSh"/"File" ---
$dup.TextFrame.TextRange.Text = "ZZZZ"
$dup.TextFrame.TextRange.Text = "ScSh`rFile"
$dup.TextFrame.TextRange.Font.Size = 15

$dup.Left = 193.91928863858266
$dup.Top = 241.77071386141733
$dup.Width = 48.52291298582677
$dup.Height = 43.62188916377953
